# Updates the cryptos price/volume table with refreshed values scraped from coinranking.com.
# Rows 47 and 48 also swap their Coin/Link/Price/Volume content (Monero <-> Stellar reorder).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "'" + '60.749.67'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.55%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "'" + '2.979.59'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.85%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.13%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "'" + '525.76'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.23%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "'" + '130.27'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.99%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.14%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "'" + '2.976.98'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.53%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "'" + '0.489'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.35%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "'" + '0.148'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.18%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "'" + '6.10'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.34%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "'" + '0.439'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.68%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "'" + '0.0000218'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.82%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "'" + '33.13'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -3.71%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "'" + '3.455.26'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.19%  '

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.22%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "'" + '60.787.05'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -3.53%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "'" + '2.976.86'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.23%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "'" + '6.48'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.16%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "'" + '456.72'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -5.63%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "'" + '13.04'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.93%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "'" + '0.666'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -4.37%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "'" + '6.79'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -5.02%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "'" + '78.10'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.45%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "'" + '11.76'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.65%  '

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.18%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "'" + '2.63'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.95%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "'" + '7.62'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -6.51%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "'" + '0.999'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.10%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.15%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "'" + '25.15'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.15%  '

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.71%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "'" + '54.82'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.99%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "'" + '2.24'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -6.57%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "'" + '5.27'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.48%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "'" + '5.74'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.66%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "'" + '453.32'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.49%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "'" + '3.176.52'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.29%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "'" + '0.0773'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.75%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "'" + '0.0377'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.18%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "'" + '0.116'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.30%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "'" + '8.00'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.94%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "'" + '2.40'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -8.43%  '

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.05%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "'" + '0.242'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.62%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "'" + '25.32'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.61%  '

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "'" + 'Monero'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "'" + 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "'" + '118.05'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.41%  '

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "'" + 'Stellar'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "'" + 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "'" + '0.108'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.01%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "'" + '1.94'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -4.77%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "'" + '0.0₃0493'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -9.47%  '

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +9.72%  '
